$d = $word.ActiveDocument

# Locate the paragraph ending with "...今天天气不错" (June 7th entry)
$found = $d.Content.Find.Execute("晴，今天是高考第一天，上午考语文，下午考数学。今天天气不错", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$target = $d.Content.Find.Parent
$target.Collapse(0)
$target.InsertAfter("，心情也很好。")

$d.Save()
